# Green Hydrogen Electricity Supply Shareweights - update offshore wind row
#
# The "offshore wind" row (row 15) in the GHESS sheet has its Max Fraction
# of Production (dimensionless) shareweights changed from 1 to 0 for every
# forecast year column (B:AE, i.e. 2021-2050).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GHESS")

# Offshore wind shareweights (row 15) -> 0 for all years
$ws.Range("B15:AE15").Value = 0

# Reproduce the author's on-screen selection left behind in the saved file
# (the whole updated offshore-wind block, anchored at B14).
$ws.Activate()
$ws.Range("B14:AE15").Select()

# Restore "About" as the active sheet/tab, matching the saved workbook view.
$wb.Worksheets.Item("About").Activate()
